# Update the cached "datetimeFigureOut" date-placeholder text that PowerPoint
# stamps into every slide layout, the slide master, and the notes master
# (Insert > Header & Footer > Date and time). The deck was last saved while
# the field showed "2/14/2022"; re-saving it later re-cached the field text
# as "3/21/2022" everywhere it appears.

function Set-DatePlaceholderText($container, [string]$oldText, [string]$newText) {
    $count = $container.Shapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.HasTextFrame -eq -1) {
            if ($sh.TextFrame.TextRange.Text -eq $oldText) {
                $sh.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$p = $ppt.ActivePresentation
$oldDate = "2/14/2022"
$newDate = "3/21/2022"

# Slide master's own date placeholder.
Set-DatePlaceholderText $p.SlideMaster $oldDate $newDate

# Every slide layout's date placeholder.
$master = $p.SlideMaster
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    Set-DatePlaceholderText $master.CustomLayouts.Item($L) $oldDate $newDate
}

# Notes master's date placeholder (writing straight to the shape's text
# does not stick for the notes master in this host, so go through the
# HeadersFooters.DateAndTime surface instead, which does persist).
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate
